$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Modelo"
$ws.Range("F2").Value = "Pipeline(steps=[('model', GradientBoostingRegressor(n_estimators=150))])"

# Copy formatting from E1 (existing header cell) to F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
